# WebForm User Assignment execution
# Updates the randomly-assigned phone numbers (column F / PN_Value) and the
# Match1UserPos / Match2UserPos computed columns (AM2 / AN2) produced by a
# fresh run of the WebForm "OneYN_TwoYN" automation test.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$addr,
        [string]$text
    )
    $rng = $ws.Range($addr)
    # Force the cell to be interpreted/stored as text (shared string) rather
    # than being auto-coerced into a number, then drop the temporary number
    # format so the cell keeps using the default (general) style - matching
    # how the rest of the sheet's text cells are stored.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Column F (PN_Value) - new phone numbers assigned for this run
Set-TextValue "F2" "9840039077"
Set-TextValue "F3" "9840089854"
Set-TextValue "F4" "9840016875"
Set-TextValue "F5" "9840023413"
Set-TextValue "F6" "9840099492"
Set-TextValue "F7" "9840009690"
Set-TextValue "F8" "9840048061"
Set-TextValue "F9" "9840059353"
Set-TextValue "F10" "9840027940"

# Match1UserPos / Match2UserPos recalculated for this run
Set-TextValue "AM2" "0"
Set-TextValue "AN2" "2"

# Reflect where the user ended up scrolled/selected in the sheet afterwards
$ws.Range("Y2").Select()
